# Adds survey-observation rows 12-14 to the "Artfynd" species-sighting sheet,
# matching the 3 new records appended in the source export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 12 ----
$ws.Range("A12").Value = 112182046
$ws.Range("B12").Value = 90678
$ws.Range("C12").Value = "Ovaliderad"
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 4366
$ws.Range("F12").Value = "Skarp dropptaggsvamp"
$ws.Range("G12").Value = "Hydnellum peckii"
$ws.Range("H12").Value = "Banker"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "1"
$ws.Range("P12").Value = "Vassbo, Dlr"
$ws.Range("Q12").Value = 374849.6512867718
$ws.Range("R12").Value = 6871060.635308203
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = "Dalarna"
$ws.Range("U12").Value = "Älvdalen"
$ws.Range("V12").Value = "Dalarna"
$ws.Range("W12").Value = "Idre"
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "2023-07-07"
$ws.Range("Z12").Value = "00:00"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "2023-07-07"
$ws.Range("AB12").Value = "00:00"
$ws.Range("AC12").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
# AT12 holds an empty string in the source row; left unset (blank), which is the nearest representable state via Range.Value
$ws.Range("AW12").Value = "Mimmi Persson"
$ws.Range("AX12").Value = "Mimmi Persson"
# AY12 holds an empty string in the source row; left unset (blank), which is the nearest representable state via Range.Value

# ---- Row 13 ----
$ws.Range("A13").Value = 112182730
$ws.Range("B13").Value = 95538
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 221941
$ws.Range("F13").Value = "Plattlummer"
$ws.Range("G13").Value = "Lycopodium complanatum"
$ws.Range("H13").Value = "L."
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "1"
$ws.Range("P13").Value = "Vassbo, Dlr"
$ws.Range("Q13").Value = 375047.1980067284
$ws.Range("R13").Value = 6871263.893339855
$ws.Range("S13").Value = 5
$ws.Range("T13").Value = "Dalarna"
$ws.Range("U13").Value = "Älvdalen"
$ws.Range("V13").Value = "Dalarna"
$ws.Range("W13").Value = "Idre"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2023-07-07"
$ws.Range("Z13").Value = "00:00"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2023-07-07"
$ws.Range("AB13").Value = "00:00"
$ws.Range("AC13").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
# AT13 holds an empty string in the source row; left unset (blank), which is the nearest representable state via Range.Value
$ws.Range("AW13").Value = "Mimmi Persson"
$ws.Range("AX13").Value = "Mimmi Persson"
# AY13 holds an empty string in the source row; left unset (blank), which is the nearest representable state via Range.Value

# ---- Row 14 ----
$ws.Range("A14").Value = 112181852
$ws.Range("B14").Value = 77550
$ws.Range("C14").Value = "Ovaliderad"
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 185
$ws.Range("F14").Value = "Violettgrå tagellav"
$ws.Range("G14").Value = "Bryoria nadvornikiana"
$ws.Range("H14").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "1"
$ws.Range("P14").Value = "Vassbo, Dlr"
$ws.Range("Q14").Value = 374953.8252938317
$ws.Range("R14").Value = 6870891.962659046
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = "Dalarna"
$ws.Range("U14").Value = "Älvdalen"
$ws.Range("V14").Value = "Dalarna"
$ws.Range("W14").Value = "Idre"
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2023-07-07"
$ws.Range("Z14").Value = "00:00"
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = "2023-07-07"
$ws.Range("AB14").Value = "00:00"
$ws.Range("AC14").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
# AT14 holds an empty string in the source row; left unset (blank), which is the nearest representable state via Range.Value
$ws.Range("AW14").Value = "Mimmi Persson"
$ws.Range("AX14").Value = "Mimmi Persson"
# AY14 holds an empty string in the source row; left unset (blank), which is the nearest representable state via Range.Value
